$data = @(
  @("BNK", "2024-06-14", "에스오에스랩", "한국", "한국, BNK", "2024-06-19", "2024-06-25", 1150, 2000000, 11500, 0, 5),
  @("DB", "2024-06-05", "디비금융스팩12호", "DB", "DB", "2024-06-11", "2024-06-18", 10000, 5000000, 2000, 0, 100),
  @("IBK", "2024-06-10", "한중엔시에스", "IBK", "IBK", "2024-06-13", "2024-06-24", 48000, 1600000, 30000, 0, 100),
  @("KB", "2024-05-07", "KB제28호스팩", "KB", "KB", "2024-05-10", "2024-05-17", 10000, 5000000, 2000, 0, 100),
  @("KB", "2024-06-11", "KB제29호스팩", "KB", "KB", "2024-06-14", "2024-06-21", 12000, 6000000, 2000, 0, 100),
  @("KB", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 215255.4, 8900000, 83400, 0, 29),
  @("NH", "2024-05-07", "아이씨티케이", "NH", "NH", "2024-05-10", "2024-05-17", 39400, 1970000, 20000, 0, 100),
  @("대신", "2024-06-05", "라메디텍", "대신", "대신", "2024-06-11", "2024-06-17", 20768, 1298000, 16000, 0, 100),
  @("대신", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 18556.5, 8900000, 83400, 0, 2.5),
  @("미래", "2024-06-10", "미래에셋비전스팩5호", "미래", "미래", "2024-06-13", "2024-06-19", 9500, 4750000, 2000, 0, 100),
  @("미래", "2024-05-20", "미래에셋비전스팩4호", "미래", "미래", "2024-05-23", "2024-05-29", 13300, 6650000, 2000, 0, 100),
  @("미래", "2024-06-13", "미래에셋비전스팩6호", "미래", "미래", "2024-06-18", "2024-06-24", 12900, 6450000, 2000, 0, 100),
  @("삼성", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 18556.5, 8900000, 83400, 0, 2.5),
  @("삼성", "2024-06-03", "그리드위즈", "삼성", "삼성", "2024-06-07", "2024-06-14", 56000, 1400000, 40000, 0, 100),
  @("삼성", "2024-05-13", "노브랜드", "삼성", "삼성", "2024-05-17", "2024-05-23", 16800, 1200000, 14000, 0, 100),
  @("신한", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 74226, 8900000, 83400, 0, 10),
  @("유비에스", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 170719.8, 8900000, 83400, 0, 23),
  @("제이피모간회사", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 170719.8, 8900000, 83400, 0, 23),
  @("하나", "2024-04-25", "HD현대마린솔루션", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점", "KB, 유비에스리미티드(영업소), 제이피모간회사 서울지점, 신한, 하나, 대신, 삼성", "2024-04-30", "2024-05-08", 74226, 8900000, 83400, 0, 10),
  @("한국", "2024-06-18", "하이젠알앤엠", "한국", "한국", "2024-06-21", "2024-06-27", 23800, 3400000, 7000, 0, 100),
  @("한국", "2024-06-10", "한국제14호스팩", "한국", "한국", "2024-06-13", "2024-06-19", 8000, 4000000, 2000, 0, 100),
  @("한국", "2024-06-14", "에스오에스랩", "한국", "한국, BNK", "2024-06-19", "2024-06-25", 21850, 2000000, 11500, 0, 95),
  @("한국", "2024-06-17", "한국제15호스팩", "한국", "한국", "2024-06-20", "2024-06-26", 12500, 6250000, 2000, 0, 100),
  @("한국", "2024-06-10", "씨어스테크놀로지", "한국", "한국", "2024-06-13", "2024-06-19", 22100, 1300000, 17000, 0, 100),
  @("한국", "2024-04-24", "코칩", "한국", "한국", "2024-04-29", "2024-05-07", 27000, 1500000, 18000, 0, 100),
  @("현대차", "2024-06-11", "에이치엠씨제7호스팩", "현대차", "현대차", "2024-06-14", "2024-06-24", 14000, 7000000, 2000, 0, 100)
)
$ws = $excel.ActiveWorkbook.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = "'" + $row[5]
    $ws.Cells.Item($r, 7).Value  = "'" + $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}

# Strip the "stored as text" quote-prefix styling picked up from the
# apostrophe-forced text entry in the date-like columns, restoring the
# default (un-styled) look used throughout the data rows.
$ws.Range("A2:L27").Style = "Normal"
